$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the active selection from E22 to D14
$ws.Range("D14").Select() | Out-Null

# Unhide the previously-hidden rows 11, 16, 23, 24, 25, 30, 37
$ws.Rows.Item(11).Hidden = $false
$ws.Rows.Item(16).Hidden = $false
$ws.Rows.Item(23).Hidden = $false
$ws.Rows.Item(24).Hidden = $false
$ws.Rows.Item(25).Hidden = $false
$ws.Rows.Item(30).Hidden = $false
$ws.Rows.Item(37).Hidden = $false

# Remove the autofilter criteria on column E (field 5), keeping column D's filter intact
$ws.AutoFilter.Range.AutoFilter(5) | Out-Null
